$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99, shifting existing rows 99:197 down to 100:198
$ws.Rows("99:99").Insert()

# Populate the newly inserted row 99 with the new data record
$ws.Range("A99").Value = 11
$ws.Range("B99").Value = "Vega Monumental Concepción"
$ws.Range("C99").Value = "Bíobío"
$ws.Range("D99").Value = 45225
$ws.Range("E99").Value = 8
$ws.Range("F99").Value = 100112001
$ws.Range("G99").Value = "Berenjena"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 100
$ws.Range("K99").Value = 10000
$ws.Range("L99").Value = 10000
$ws.Range("M99").Value = 10000
$ws.Range("N99").Value = "`$/caja 50 unidades"
$ws.Range("O99").Value = "Región de Arica y Parinacota"
$ws.Range("P99").Value = 200
$ws.Range("Q99").Value = 50
$ws.Range("R99").Value = "Hortaliza"
